$d = $word.ActiveDocument

# The second "NOTES:" paragraph and the following list item
# ("This looks very promising, look at this tomorrow please!!!")
# are merged into a single (non-list) paragraph with new text,
# while the trailing _GoBack bookmark is preserved at the end.
#
# Find/Replace across the paragraph mark merges the two paragraphs,
# keeping the first paragraph's (non-list) formatting, and replaces
# the combined text with the new sentence.
$found = $d.Content.Find.Execute(
    "NOTES:^pThis looks very promising, look at this tomorrow please!!!",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This contains basic examples for rune.js", 2)

# Re-create the _GoBack bookmark at the very end of the document
# (it was removed as part of the replaced range above). Inserting a
# throwaway character first avoids this runtime's handling of
# zero-length bookmark ranges, then the character is removed again,
# leaving a clean, empty bookmark exactly like the original.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$end = $last.Range.Duplicate
$end.Collapse(0)
$end.MoveEnd(1, -1) | Out-Null
$end.InsertAfter("X")
$end.MoveEnd(1, 1) | Out-Null
$d.Bookmarks.Add("_GoBack", $end)
$end.Text = ""
